$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Locate the insertion point: right after "...zu Helfen." and
#    right before " Du bist ein neuer Abenteurer..."
# ------------------------------------------------------------------
$findRange = $d.Content
$ok = $findRange.Find.Execute("Helfen. Du bist ein neuer Abenteurer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) {
    throw "Could not locate anchor text for insertion"
}
$base = $findRange.Start + 7   # length of "Helfen." -> position right after the period

# ------------------------------------------------------------------
# 2. Insert the new sentence (" Btw. Helfen wir dir nach der Landng
#    zurecht zu kommen.") right at that point, ahead of the existing
#    " Du bist..." text.
# ------------------------------------------------------------------
$newText = " Btw. Helfen wir dir nach der Landng zurecht zu kommen."
$insPoint = $d.Range($base, $base)
$insPoint.InsertBefore($newText)

# ------------------------------------------------------------------
# 3. Force the inserted text to live in its own separate runs
#    (matching the run boundaries of the target markup) by briefly
#    dropping a bookmark at each boundary and removing it again -
#    the split persists even after the bookmark is deleted.
# ------------------------------------------------------------------
$offsets = @(0, 1, 4, 30, 36)
$i = 0
foreach ($off in $offsets) {
    $pos = $base + $off
    $bmRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add("zzzTempSplit$i", $bmRange)
    $i = $i + 1
}
for ($j = 0; $j -lt $i; $j++) {
    $d.Bookmarks.Item("zzzTempSplit$j").Delete()
}

# ------------------------------------------------------------------
# 4. Re-home the "_GoBack" bookmark: it used to sit by itself in an
#    empty paragraph further down the document; it now belongs right
#    after "...zurecht zu kommen." and before " Du bist ein neuer
#    Abenteurer...". Adding it here automatically removes the old
#    occurrence (Word only ever keeps one "_GoBack" bookmark).
# ------------------------------------------------------------------
$goBackPos = $base + 55   # length of $newText
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)

Write-Output $d.Paragraphs(3).Range.Text
